$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.319.05'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '1.855.84'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '239.34'
$ws.Range("E5").Value = '  -1.28%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = '0.6966'
$ws.Range("E6").Value = '  -6.41%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '0.3082'
$ws.Range("E8").Value = '  -2.19%  '
$ws.Range("D9").Value = '0.07554'
$ws.Range("E9").Value = '  +4.33%  '
$ws.Range("D10").Value = '23.78'
$ws.Range("E10").Value = '  -3.80%  '
$ws.Range("D11").Value = '0.08125'
$ws.Range("E11").Value = '  -3.35%  '
$ws.Range("D12").Value = '1.873.27'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '0.7272'
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").Value = '5.208'
$ws.Range("E14").Value = '  -4.06%  '
$ws.Range("D15").Value = '89.54'
$ws.Range("E15").Value = '  -3.23%  '
$ws.Range("D16").Value = '29.449.39'
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.910'
$ws.Range("E17").Value = '  -2.96%  '
$ws.Range("D18").Value = '242.95'
$ws.Range("E18").Value = '  -1.98%  '
$ws.Range("D19").Value = '0.000007768'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("E20").Value = '  -3.21%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.126.63'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '7.629'
$ws.Range("E24").Value = '  -5.20%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.065'
$ws.Range("E25").Value = '  -2.28%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1472'
$ws.Range("E26").Value = '  -5.72%  '
$ws.Range("D27").Value = '162.06'
$ws.Range("E27").Value = '  -1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.940'
$ws.Range("E29").Value = '  -4.76%  '
$ws.Range("D30").Value = '1.405'
$ws.Range("E30").Value = '  -7.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.510'
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("D32").Value = '4.427'
$ws.Range("E32").Value = '  -3.79%  '
$ws.Range("D33").Value = '4.053'
$ws.Range("E33").Value = '  -5.36%  '
$ws.Range("D34").Value = '0.05231'
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("D35").Value = '1.199'
$ws.Range("E35").Value = '  -3.16%  '
$ws.Range("D36").Value = '0.7195'
$ws.Range("E36").Value = '  -4.25%  '
$ws.Range("D37").Value = '0.9993'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '2.664'
$ws.Range("E38").Value = '  -0.99%  '
$ws.Range("D39").Value = '0.01868'
$ws.Range("E39").Value = '  -5.04%  '
$ws.Range("D40").Value = '2.713'
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("D41").Value = '0.8845'
$ws.Range("E41").Value = '  +3.38%  '
$ws.Range("D42").Value = '0.4312'
$ws.Range("E42").Value = '  -5.31%  '
$ws.Range("D43").Value = '5.899'
$ws.Range("E43").Value = '  -2.63%  '
$ws.Range("D44").Value = '70.13'
$ws.Range("E44").Value = '  -3.31%  '
$ws.Range("D45").Value = '1.046.13'
$ws.Range("E45").Value = '  -6.04%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.60'
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").Value = '7.272'
$ws.Range("E48").Value = '  -4.61%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.752'
$ws.Range("E49").Value = '  -5.77%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.016.45'
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").Value = '9.273'
$ws.Range("E51").Value = '  -2.32%  '
